$d = $word.ActiveDocument

# --- Locate the "Meeting described in the tasks." paragraph and the
#     "---End of week 6---" paragraph that immediately follows it ---
$meetingPara = $null
$endWeek6Para = $null
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Meeting described in the tasks.*") {
        $meetingPara = $p
    }
    if ($meetingPara -ne $null -and $endWeek6Para -eq $null -and $p.Range.Text -eq "---End of week 6---`r") {
        $endWeek6Para = $p
    }
}

# --- Remove the old "_GoBack" bookmark / paragraph break between them so the
#     "---End of week 6---" paragraph (and everything after it) can be
#     rebuilt fresh, matching what a human re-typing this section would
#     produce ---
$joinStart = $meetingPara.Range.End - 1
$joinEnd = $endWeek6Para.Range.End
$joinRange = $d.Range($joinStart, $joinEnd)
$joinRange.Delete()

# --- Re-fetch the (now merged) "Meeting described..." paragraph and use it
#     as the anchor for four freshly-created paragraphs ---
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Meeting described in the tasks.*") {
        $anchorPara = $p
    }
}

$insertPoint = $anchorPara.Range
$insertPoint.Collapse(0)   # wdCollapseEnd

# 1) "---End of week 6---"
$insertPoint.InsertParagraphAfter()
$p1 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $insertPoint.Start -and $p.Range.Text -eq "`r") {
        $p1 = $p
    }
}
$p1Text = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$p1Text.Text = "---End of week 6---"

# 2) "---Week 7---"
$insertPoint = $d.Range($p1.Range.Start, $p1.Range.Start)
$insertPoint.Collapse(0)
$insertPoint.MoveEnd(1, 1) | Out-Null
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $insertPoint.Start -and $p.Range.Text -eq "`r") {
        $p2 = $p
    }
}
$p2Text = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$p2Text.Text = "---Week 7---"

# 3) "Started with a demo project ... had a lot of re-working to do."
#    (with the "_GoBack" bookmark placed right after "...project and ")
$insertPoint = $d.Range($p2.Range.Start, $p2.Range.Start)
$insertPoint.Collapse(0)
$insertPoint.MoveEnd(1, 1) | Out-Null
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()
$p3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $insertPoint.Start -and $p.Range.Text -eq "`r") {
        $p3 = $p
    }
}
$p3Start = $p3.Range.Start
$beforeBookmark = "Started with a demo project to understand how layers may work. Once I did the layers demo, I went to the white board and drew up a basic design for layers. I implemented layers into my actual project and "
$afterBookmark = "had a lot of re-working to do."
$p3Text = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$p3Text.Text = $beforeBookmark + $afterBookmark

$bmPos = $p3Start + $beforeBookmark.Length
$bmRange = $d.Range($bmPos, $bmPos)
$bmRange.Bookmarks.Add("_GoBack")

# 4) "---End of week 7---"
$insertPoint = $d.Range($p3Start, $p3Start)
$insertPoint.Collapse(0)
$insertPoint.MoveEndUntil("`r", 1000000) | Out-Null
$insertPoint.Collapse(0)
$insertPoint.MoveEnd(1, 1) | Out-Null
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()
$p4 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $insertPoint.Start -and $p.Range.Text -eq "`r") {
        $p4 = $p
    }
}
$p4Text = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$p4Text.Text = "---End of week 7---"

Write-Host "Done"
